$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Solder and Soldering Iron Kit rows as "Owned" (checkmark instead of X)
$ws.Range("L13").Value = "√"
$ws.Range("L14").Value = "√"

# Extend "Total Owned" sum to include rows 13 and 14 (now owned)
$ws.Range("F17").Formula = "=SUM(K3:K14)"

# "Total Unowned" no longer has any unowned rows -> hardcoded 0
$ws.Range("F18").Value = 0

# Update the view: scroll so B4 is the top-left visible cell, and select G21
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G21").Select()
